# Updates the "Corte" sheet with a new clip to cut: the intro vignette
# (vinheta de inicio) for the "AlgoJS - Algoritmos" course, replacing
# the previous sample row, and adds the formula/TXT columns that build
# the pipe-delimited cut-parameters string consumed by the external
# tooling/readme.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Corte")

$ws.Range("C2").Value = "M-AlgoJS - Algoritmos.mp4"
$ws.Range("B2").Value = "00:04.00"
$ws.Range("A2").Value = "./arquivos/sem_edicao/Aula01.mp4"

$ws.Range("D2").Formula = "=CONCATENATE(A2,"";"",B2,"";"",C2)"
$ws.Range("E2").Value = "./arquivos/sem_edicao/Aula01.mp4;00:04.00;M-AlgoJS - Algoritmos.mp4"

$ws.Range("E2").Select()
